$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '39.108.57'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -3.10%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.187.24'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -7.54%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '296.52'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -4.50%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '81.41'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -5.40%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -4.59%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -5.15%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -7.63%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '28.98'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -4.53%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.22'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -10.47%  '

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.63%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.23'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.56%  '

$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.529.32'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -7.60%  '

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -7.15%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.200.45'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -7.25%  '

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -6.58%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '38.991.15'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -3.39%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -4.58%  '

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -7.09%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '64.76'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -5.16%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -5.62%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '225.10'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -4.14%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.40'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -6.86%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.05%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.47'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -5.25%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.54%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.38%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '148.92'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.53%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '31.49'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -8.87%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.13%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -7.97%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.99%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0689'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -5.56%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -4.50%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0964'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.65%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.16'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.34%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -6.44%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -5.53%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.59'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -6.42%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.888.74'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -4.03%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.08'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -12.59%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -3.87%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.91'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -5.06%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '15.94'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -9.88%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -3.25%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.401.77'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -7.50%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '70.88'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.69%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '86.52'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -6.90%  '
